# Adds a fourth "Commit 4" results block (rows 70-85) to Sheet1, mirroring
# the existing "Commit 3" block (rows 52-67) with updated numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New shared-string label for the section header.
# ---------------------------------------------------------------------
$ws.Range("A70").Value = "Commit 4"

# ---------------------------------------------------------------------
# 2) "MARS Tool Output" / "Calulations" banner row (row 71, bold 14pt).
#    Copy formatting from the equivalent row in the "Commit 3" block
#    (row 53) so we reuse the existing cell styles instead of minting
#    new ones.
# ---------------------------------------------------------------------
$ws.Range("A71").Value = "MARS Tool Output"
$ws.Range("D71").Value = "Calulations"

$ws.Range("A53").Copy() | Out-Null
$ws.Range("A71").PasteSpecial(-4122) | Out-Null
$ws.Range("B53").Copy() | Out-Null
$ws.Range("B71").PasteSpecial(-4122) | Out-Null
$ws.Range("C53").Copy() | Out-Null
$ws.Range("C71").PasteSpecial(-4122) | Out-Null
$ws.Range("D53").Copy() | Out-Null
$ws.Range("D71").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(71).RowHeight = 18.5

# ---------------------------------------------------------------------
# 3) "Instruction Statistics Tool" sub-table (rows 73-79).
# ---------------------------------------------------------------------
$ws.Range("A73").Value = "Instruction Statistics Tool"
$ws.Range("A55").Copy() | Out-Null
$ws.Range("A73").PasteSpecial(-4122) | Out-Null

$ws.Range("A74").Value = "Instruction type"
$ws.Range("B74").Value = "Count"
$ws.Range("D74").Value = "Adjusted count"
$ws.Range("E74").Value = "CPI"
$ws.Range("F74").Value = "Total cycles"

$ws.Range("A56").Copy() | Out-Null
$ws.Range("A74").PasteSpecial(-4122) | Out-Null
$ws.Range("B56").Copy() | Out-Null
$ws.Range("B74").PasteSpecial(-4122) | Out-Null
$ws.Range("D56").Copy() | Out-Null
$ws.Range("D74").PasteSpecial(-4122) | Out-Null
$ws.Range("E56").Copy() | Out-Null
$ws.Range("E74").PasteSpecial(-4122) | Out-Null
$ws.Range("F56").Copy() | Out-Null
$ws.Range("F74").PasteSpecial(-4122) | Out-Null

# ALU
$ws.Range("A75").Value = "ALU"
$ws.Range("B75").Value = 3625
$ws.Range("D75").Formula = "=B75"
$ws.Range("E75").Value = 1
$ws.Range("F75").Formula = "=D75*E75"

# Jump
$ws.Range("A76").Value = "Jump"
$ws.Range("B76").Value = 298
$ws.Range("D76").Formula = "=B76"
$ws.Range("E76").Value = 1
$ws.Range("F76").Formula = "=D76*E76"

# Branch
$ws.Range("A77").Value = "Branch"
$ws.Range("B77").Value = 965
$ws.Range("D77").Formula = "=B77"
$ws.Range("E77").Value = 2
$ws.Range("F77").Formula = "=D77*E77"

# Memory
$ws.Range("A78").Value = "Memory"
$ws.Range("B78").Value = 620

# Other
$ws.Range("A79").Value = "Other"
$ws.Range("B79").Value = 761
$ws.Range("D79").Formula = "=B79-(B83+B84-B78)"
$ws.Range("E79").Value = 5
$ws.Range("F79").Formula = "=D79*E79"

# ---------------------------------------------------------------------
# 4) "Data Cache Simulation Tool" sub-table (rows 81-84).
# ---------------------------------------------------------------------
$ws.Range("A81").Value = "Data Cache Simulation Tool"
$ws.Range("A63").Copy() | Out-Null
$ws.Range("A81").PasteSpecial(-4122) | Out-Null

$ws.Range("A82").Value = "Access"
$ws.Range("B82").Value = "Count"
$ws.Range("A64").Copy() | Out-Null
$ws.Range("A82").PasteSpecial(-4122) | Out-Null
$ws.Range("B64").Copy() | Out-Null
$ws.Range("B82").PasteSpecial(-4122) | Out-Null

# Cache hit
$ws.Range("A83").Value = "Cache hit"
$ws.Range("B83").Value = 361
$ws.Range("D83").Formula = "=B83"
$ws.Range("E83").Value = 2
$ws.Range("F83").Formula = "=D83*E83"

# Cache miss
$ws.Range("A84").Value = "Cache miss"
$ws.Range("B84").Value = 351
$ws.Range("D84").Formula = "=B84"
$ws.Range("E84").Value = 40
$ws.Range("F84").Formula = "=D84*E84"

# ---------------------------------------------------------------------
# 5) Totals row (row 85, bold 14pt + top border).
# ---------------------------------------------------------------------
$ws.Range("F85").Formula = "=SUM(F75:F84)"
$ws.Range("F67").Copy() | Out-Null
$ws.Range("F85").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(85).RowHeight = 18.5

# ---------------------------------------------------------------------
# 6) Update the view to mirror the scroll/selection described in the diff.
# ---------------------------------------------------------------------
$ws.Range("B76").Select()
